{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list to use punchy,\n// impact-focused accomplishment statements (vs. the old job-duty style\n// bullets), and drop the two bullets that duplicated Core-Competency-style\n// tooling callouts \u2014 matching the target diff exactly.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Some of this section's OLD bullet text (e.g. the \"Trigonometric\n// algorithm...\" line) is duplicated verbatim elsewhere in the resume (the\n// Siege Analytics job bullets), so a document-wide text match would be\n// ambiguous. Find the \"KEY ACHIEVEMENTS AND IMPACT\" heading first and only\n// touch paragraphs between it and the next heading (\"TECHNICAL SKILLS\").\nlet sectionStart = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionStart = i;\n    break;\n  }\n}\nif (sectionStart === -1) {\n  throw new Error(\"KEY ACHIEVEMENTS AND IMPACT section not found\");\n}\n\nlet sectionEnd = items.length;\nfor (let i = sectionStart + 1; i < items.length; i++) {\n  if (items[i].text === \"TECHNICAL SKILLS\") {\n    sectionEnd = i;\n    break;\n  }\n}\n\n// Exact OLD bullet text -> NEW replacement text, scoped to this section.\nconst replacements = new Map([\n  [\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  ],\n  [\n    \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    \"\u2022 $4.7M savings enabled nonprofit access\",\n  ],\n  [\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  ],\n  [\n    \"\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\",\n    \"\u2022 178% accuracy improvement in racial classification algorithms\",\n  ],\n]);\n\n// Paragraphs to remove outright (within the section).\nconst toDelete = new Set([\n  \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n]);\n\nconst deleteTargets = [];\nfor (let i = sectionStart; i < sectionEnd; i++) {\n  const text = items[i].text;\n  if (replacements.has(text)) {\n    items[i].insertText(replacements.get(text), \"Replace\");\n  } else if (toDelete.has(text)) {\n    deleteTargets.push(items[i]);\n  }\n}\n\nfor (const p of deleteTargets) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list to use punchy,\n# impact-focused accomplishment statements (vs. the old job-duty style\n# bullets), and drop the two bullets that duplicated Core-Competency-style\n# tooling callouts - matching the target diff exactly.\n\n$d = $word.ActiveDocument\n\n# --- Locate the \"KEY ACHIEVEMENTS AND IMPACT\" section -----------------\n# Some bullet text in this section (e.g. the \"Trigonometric algorithm...\"\n# line) is duplicated verbatim elsewhere in the resume (Siege Analytics job\n# bullets), so a plain document-wide Find/Replace would be ambiguous. Scope\n# all work to the paragraph range between the \"KEY ACHIEVEMENTS AND IMPACT\"\n# heading and the next heading (\"TECHNICAL SKILLS\").\n$count = $d.Paragraphs.Count\n$sectionStart = -1\n$sectionEnd = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($sectionStart -eq -1 -and $txt -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionStart = $i\n    } elseif ($sectionStart -ne -1 -and $sectionEnd -eq -1 -and $txt -eq \"TECHNICAL SKILLS\") {\n        $sectionEnd = $i\n    }\n}\n\nif ($sectionStart -eq -1) {\n    throw \"KEY ACHIEVEMENTS AND IMPACT section not found\"\n}\nif ($sectionEnd -eq -1) {\n    $sectionEnd = $count + 1\n}\n\n$bullet = [char]0x2022\n\n# Exact OLD bullet text -> NEW replacement text, scoped to this section only.\n$replacements = @{}\n$replacements[$bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"] = $bullet + \" Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n$replacements[$bullet + \" Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\"] = $bullet + \" `$4.7M savings enabled nonprofit access\"\n$replacements[$bullet + \" Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\"] = $bullet + \" Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n$replacements[$bullet + \" Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\"] = $bullet + \" 178% accuracy improvement in racial classification algorithms\"\n\n# Paragraphs to remove outright (within the section).\n$toDelete = @(\n    ($bullet + \" Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\"),\n    ($bullet + \" Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\")\n)\n\n# Walk the section range backwards so paragraph indices for not-yet-visited\n# items stay valid as deletions shift the collection.\nfor ($i = $sectionEnd - 1; $i -ge $sectionStart; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13)\n\n    if ($replacements.ContainsKey($txt)) {\n        $p.Range.Text = $replacements[$txt]\n    } elseif ($toDelete -contains $txt) {\n        $p.Range.Delete()\n    }\n}\n"}
